$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute("`${TienAn} `${TienSu}", $true, $false, $false, $false, $false,
              $true, 1, $false, "", 2)
